$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at spreadsheet row 22, pushing the existing
# rows 22..147 down to 23..148 (weekly price update prepended to the
# "Ajo" / Feria Lagunitas de Puerto Montt series).
$ws.Rows.Item(22).Insert()

$ws.Range("A22").Value = 4
$ws.Range("B22").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C22").Value = 'Los Lagos'
$ws.Range("D22").Value = 44462
$ws.Range("E22").Value = 10
$ws.Range("F22").Value = 100112003
$ws.Range("G22").Value = 'Ajo'
$ws.Range("H22").Value = 'Chino'
$ws.Range("I22").Value = 'Primera'
$ws.Range("J22").Value = 100
$ws.Range("K22").Value = 17000
$ws.Range("L22").Value = 17000
$ws.Range("M22").Value = 17000
$ws.Range("N22").Value = '$/caja 10 kilos'
$ws.Range("O22").Value = 'China'
$ws.Range("P22").Value = 1700
$ws.Range("Q22").Value = 10
$ws.Range("R22").Value = 'Hortaliza'
